$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.284.39'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '1.849.46'
$ws.Range("E3").Value = '  +0.83%  '
$ws.Range("D4").Value = '''1.004'
$ws.Range("E4").Value = '  -0.37%  '
$ws.Range("D5").Value = '''313.81'
$ws.Range("E5").Value = '  +0.37%  '
$ws.Range("D6").Value = '''1.003'
$ws.Range("E6").Value = '  -0.25%  '
$ws.Range("D7").Value = '''0.4606'
$ws.Range("E7").Value = '  -1.02%  '
$ws.Range("E8").Value = '  -0.27%  '
$ws.Range("D9").Value = '''0.07285'
$ws.Range("E9").Value = '  -1.02%  '
$ws.Range("D10").Value = '''0.8837'
$ws.Range("E10").Value = '  +0.90%  '
$ws.Range("D11").Value = '''19.88'
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").Value = '''0.07799'
$ws.Range("E12").Value = '  -0.03%  '
$ws.Range("D13").Value = '1.864.38'
$ws.Range("E13").Value = '  +0.16%  '
$ws.Range("E14").Value = '  +0.30%  '
$ws.Range("D15").Value = '''6.540'
$ws.Range("E15").Value = '  -0.54%  '
$ws.Range("D16").Value = '''91.52'
$ws.Range("E16").Value = '  -0.61%  '
$ws.Range("D17").Value = '''1.004'
$ws.Range("E17").Value = '  -0.33%  '
$ws.Range("D18").Value = '''0.000008943'
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("D19").Value = '''1.002'
$ws.Range("E19").Value = '  -0.21%  '
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").Value = '27.308.23'
$ws.Range("E21").Value = '  +1.86%  '
$ws.Range("D22").Value = '''5.115'
$ws.Range("E22").Value = '  -0.72%  '
$ws.Range("D23").Value = '''10.53'
$ws.Range("E23").Value = '  -0.70%  '
$ws.Range("D24").Value = '2.061.14'
$ws.Range("E24").Value = '  -3.31%  '
$ws.Range("D25").Value = '''1.929'
$ws.Range("E25").Value = '  +5.40%  '
$ws.Range("D26").Value = '''151.63'
$ws.Range("E26").Value = '  -0.40%  '
$ws.Range("D27").Value = '''18.38'
$ws.Range("E27").Value = '  -0.03%  '
$ws.Range("D28").Value = '''2.049'
$ws.Range("E28").Value = '  -2.39%  '
$ws.Range("D29").Value = '''115.88'
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = '''5.067'
$ws.Range("E30").Value = '  -0.53%  '
$ws.Range("D31").Value = '''0.08833'
$ws.Range("E31").Value = '  -0.48%  '
$ws.Range("D32").Value = '''3.119'
$ws.Range("E32").Value = '  +5.39%  '
$ws.Range("D33").Value = '''0.7685'
$ws.Range("E33").Value = '  +5.11%  '
$ws.Range("D34").Value = '''1.167'
$ws.Range("E34").Value = '  +2.28%  '
$ws.Range("D35").Value = '''4.490'
$ws.Range("E35").Value = '  +0.80%  '
$ws.Range("D36").Value = '''2.643'
$ws.Range("E36").Value = '  +5.67%  '
$ws.Range("D37").Value = '''1.080'
$ws.Range("E37").Value = '  +0.78%  '
$ws.Range("D38").Value = '''0.01955'
$ws.Range("E38").Value = '  +0.28%  '
$ws.Range("D39").Value = '''0.05231'
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("D40").Value = '''2.958'
$ws.Range("E40").Value = '  +0.75%  '
$ws.Range("D41").Value = '''7.008'
$ws.Range("E41").Value = '  -2.28%  '
$ws.Range("E42").Value = '  -1.58%  '
$ws.Range("D43").Value = '''0.1632'
$ws.Range("E43").Value = '  -0.04%  '
$ws.Range("D44").Value = '''8.373'
$ws.Range("E44").Value = '  +1.54%  '
$ws.Range("D45").Value = '''0.4799'
$ws.Range("E45").Value = '  -0.93%  '
$ws.Range("D46").Value = '''10.33'
$ws.Range("E46").Value = '  +0.55%  '
$ws.Range("E47").Value = '  -0.28%  '
$ws.Range("D48").Value = '''102.56'
$ws.Range("E48").Value = '  +0.13%  '
$ws.Range("D49").Value = '''1.647'
$ws.Range("E49").Value = '  +1.15%  '
$ws.Range("D50").Value = '''0.06222'
$ws.Range("E50").Value = '  -0.18%  '
$ws.Range("D51").Value = '''65.43'
$ws.Range("E51").Value = '  +0.93%  '
